$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1833333333333333
$ws.Range("C2").Value = 0.5666666666666667
$ws.Range("J2").Value = 0.03
$ws.Range("P2").Value = 0.11
$ws.Range("S2").Value = 0.11
$ws.Range("B3").Value = 0.005813953488372093
$ws.Range("C3").Value = 0.01162790697674419
$ws.Range("J3").Value = 0.01162790697674419
$ws.Range("P3").Value = 0.7325581395348837
$ws.Range("S3").Value = 0.2383720930232558
$ws.Range("B6").Value = 0.08860759493670886
$ws.Range("D6").Value = 0.003164556962025316
$ws.Range("E6").Value = 0.003164556962025316
$ws.Range("F6").Value = 0.08544303797468354
$ws.Range("J6").Value = 0.2626582278481013
$ws.Range("O6").Value = 0.02848101265822785
$ws.Range("Q6").Value = 0.1677215189873418
$ws.Range("R6").Value = 0.0379746835443038
$ws.Range("S6").Value = 0.3227848101265823
$ws.Range("B7").Value = 0.09025270758122744
$ws.Range("D7").Value = 0.01444043321299639
$ws.Range("F7").Value = 0.06137184115523465
$ws.Range("J7").Value = 0.09747292418772563
$ws.Range("O7").Value = 0.03249097472924187
$ws.Range("Q7").Value = 0.1877256317689531
$ws.Range("R7").Value = 0.08303249097472924
$ws.Range("S7").Value = 0.4332129963898917
$ws.Range("B8").Value = 0.07964601769911504
$ws.Range("D8").Value = 0.01238938053097345
$ws.Range("F8").Value = 0.06548672566371681
$ws.Range("J8").Value = 0.08849557522123894
$ws.Range("O8").Value = 0.03893805309734513
$ws.Range("Q8").Value = 0.1646017699115044
$ws.Range("R8").Value = 0.07610619469026549
$ws.Range("S8").Value = 0.4743362831858407
$ws.Range("B9").Value = 0.07843137254901961
$ws.Range("D9").Value = 0.0261437908496732
$ws.Range("F9").Value = 0.0718954248366013
$ws.Range("J9").Value = 0.1241830065359477
$ws.Range("O9").Value = 0.0457516339869281
$ws.Range("Q9").Value = 0.1764705882352941
$ws.Range("R9").Value = 0.1045751633986928
$ws.Range("S9").Value = 0.3725490196078431
$ws.Range("B10").Value = 0.09843400447427293
$ws.Range("D10").Value = 0.01416853094705444
$ws.Range("E10").Value = 0.001491424310216256
$ws.Range("F10").Value = 0.08053691275167785
$ws.Range("J10").Value = 0.1036539895600298
$ws.Range("O10").Value = 0.01342281879194631
$ws.Range("Q10").Value = 0.2244593586875466
$ws.Range("R10").Value = 0.06711409395973154
$ws.Range("S10").Value = 0.3967188665175242
$ws.Range("G11").Value = 0.1531322505800464
$ws.Range("J11").Value = 0.1020881670533643
$ws.Range("K11").Value = 0.2204176334106729
$ws.Range("L11").Value = 0.505800464037123
$ws.Range("S11").Value = 0.0185614849187935
$ws.Range("G12").Value = 0.7555555555555555
$ws.Range("J12").Value = 0.1644444444444444
$ws.Range("K12").Value = 0.01777777777777778
$ws.Range("L12").Value = 0.02222222222222222
$ws.Range("S12").Value = 0.04
$ws.Range("G13").Value = 0.6764705882352942
$ws.Range("J13").Value = 0.2941176470588235
$ws.Range("S13").Value = 0.02941176470588235
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.0218978102189781
$ws.Range("H15").Value = 0.2153284671532847
$ws.Range("I15").Value = 0.0291970802919708
$ws.Range("J15").Value = 0.3211678832116788
$ws.Range("K15").Value = 0.08029197080291971
$ws.Range("M15").Value = 0.0218978102189781
$ws.Range("O15").Value = 0.1021897810218978
$ws.Range("S15").Value = 0.208029197080292
$ws.Range("F16").Value = 0.02259887005649718
$ws.Range("H16").Value = 0.231638418079096
$ws.Range("I16").Value = 0.0847457627118644
$ws.Range("J16").Value = 0.3333333333333333
$ws.Range("K16").Value = 0.1186440677966102
$ws.Range("M16").Value = 0.01694915254237288
$ws.Range("N16").Value = 0.005649717514124294
$ws.Range("O16").Value = 0.04519774011299435
$ws.Range("S16").Value = 0.1412429378531073
$ws.Range("F17").Value = 0.02697495183044316
$ws.Range("H17").Value = 0.1888246628131021
$ws.Range("I17").Value = 0.06551059730250482
$ws.Range("J17").Value = 0.4007707129094412
$ws.Range("K17").Value = 0.1175337186897881
$ws.Range("M17").Value = 0.03468208092485549
$ws.Range("O17").Value = 0.07514450867052024
$ws.Range("S17").Value = 0.0905587668593449
$ws.Range("F18").Value = 0.04945054945054945
$ws.Range("H18").Value = 0.1923076923076923
$ws.Range("I18").Value = 0.06043956043956044
$ws.Range("J18").Value = 0.3791208791208791
$ws.Range("K18").Value = 0.1373626373626374
$ws.Range("M18").Value = 0.03296703296703297
$ws.Range("N18").Value = 0.005494505494505495
$ws.Range("O18").Value = 0.06043956043956044
$ws.Range("S18").Value = 0.08241758241758242
$ws.Range("F19").Value = 0.038173142467621
$ws.Range("H19").Value = 0.2283571915473756
$ws.Range("I19").Value = 0.05862304021813224
$ws.Range("J19").Value = 0.3469665985003408
$ws.Range("K19").Value = 0.1356509884117246
$ws.Range("M19").Value = 0.02385821404226312
$ws.Range("N19").Value = 0.0006816632583503749
$ws.Range("O19").Value = 0.06339468302658487
$ws.Range("S19").Value = 0.1042944785276074
